{"js": "async (context) => {\n  // The document content got reshuffled: each \"slot\" paragraph below now\n  // displays the text that used to belong to a different slot. We rewrite\n  // each paragraph's text in place (line breaks encoded as \\u000B so they\n  // serialize back out as <w:br/> between <w:t> runs, exactly like before).\n\n  const body = context.document.body;\n  const paragraphs = body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n\n  const BR = \"\\u000B\";\n\n  // Objetivos paragraph -> gets the (short) \"Programa resumido\" list.\n  paragraphs.items[5].insertText(\n    \"1. Fundamentos da Gest\u00e3o de Produ\u00e7\u00e3o\" + BR +\n    \"2. Vis\u00e3o estrat\u00e9gica da Produ\u00e7\u00e3o.\" + BR +\n    \"3. Projeto em Gest\u00e3o da Produ\u00e7\u00e3o.\" + BR +\n    \"4. Planejamento e Controle da Produ\u00e7\u00e3o\",\n    Word.InsertLocation.replace\n  );\n\n  // Docente(s) Respons\u00e1vel(eis) paragraph -> gets the old Objetivos text.\n  paragraphs.items[7].insertText(\n    \"Apresentar conceitos fundamentais de Administra\u00e7\u00e3o da Produ\u00e7\u00e3o.\",\n    Word.InsertLocation.replace\n  );\n\n  // Programa resumido paragraph -> gets the (long) \"Programa\" list.\n  paragraphs.items[9].insertText(\n    \"1 - Fundamentos da gest\u00e3o de produ\u00e7\u00e3o: modelo de transforma\u00e7\u00e3o: inputs, processo de transforma\u00e7\u00e3o e outputs. Tipos de Processo de Produ\u00e7\u00e3o\" + BR +\n    \"2 - Vis\u00e3o estrat\u00e9gica de produ\u00e7\u00e3o: Papel da fun\u00e7\u00e3o produ\u00e7\u00e3o. Objetivos de Desempenho. Estrat\u00e9gias de Produ\u00e7\u00e3o. Ciclo de Vida Produto/Servi\u00e7o.\" + BR +\n    \"3 \u2013 Projeto em Gest\u00e3o da Produ\u00e7\u00e3o: Tipos de Processos. Projeto de Produtos e Servi\u00e7os. Projeto de Rede de Opera\u00e7\u00f5es Produtivas. Arranjo F\u00edsico.\" + BR +\n    \"4 - Planejamento e Controle da Produ\u00e7\u00e3o: Material Requirement Planning (MRP), Manufacturing Resources Planning (MPRII), Enterprise Planning (ERP). Produ\u00e7\u00e3o Enxuta. Kanban. Just in Time.\",\n    Word.InsertLocation.replace\n  );\n\n  // Programa paragraph -> gets the old \"M\u00e9todo\" text.\n  paragraphs.items[11].insertText(\n    \"O desenvolvimento da disciplina ser\u00e1 baseado em leituras, aula expositiva, discuss\u00e3o e resolu\u00e7\u00e3o de estudos de caso.\",\n    Word.InsertLocation.replace\n  );\n\n  await context.sync();\n\n  // Avalia\u00e7\u00e3o paragraph keeps its bold \"M\u00e9todo:\"/\"Crit\u00e9rio:\"/\"Norma de\n  // recupera\u00e7\u00e3o:\" labels, but the values after each label shift forward.\n  // Replace back-to-front so a freshly-inserted value can never collide\n  // with (and get matched instead of) a not-yet-processed original value.\n  const p13 = paragraphs.items[13];\n\n  const normaValue = p13.search(\n    \"Prova \u00fanica com nota maior ou igual a 5,0 (cinco).\",\n    { matchCase: true }\n  );\n  normaValue.load(\"items\");\n  await context.sync();\n  normaValue.items[0].insertText(\n    \"SLACK, N. et al. Administra\u00e7\u00e3o da Produ\u00e7\u00e3o. 3 ed. S\u00e3o Paulo: Atlas, 2009.\" + BR +\n    \"CHASE, R. B. E JACOBS, F.R. Administra\u00e7\u00e3o da Produ\u00e7\u00e3o e de Opera\u00e7\u00f5es. 1 ed. Porto Alegre. Bookman. 2009.\" + BR +\n    \"CORREA, H.L.; CORREA, C.A. Administra\u00e7\u00e3o da Produ\u00e7\u00e3o e Opera\u00e7\u00f5es. 2 ed. S\u00e3o Paulo. Atlas. 2006\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n\n  const criterioValue = p13.search(\"Provas e Trabalhos\", { matchCase: true });\n  criterioValue.load(\"items\");\n  await context.sync();\n  criterioValue.items[0].insertText(\n    \"Prova \u00fanica com nota maior ou igual a 5,0 (cinco).\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n\n  const metodoValue = p13.search(\n    \"O desenvolvimento da disciplina ser\u00e1 baseado em leituras, aula expositiva, discuss\u00e3o e resolu\u00e7\u00e3o de estudos de caso.\",\n    { matchCase: true }\n  );\n  metodoValue.load(\"items\");\n  await context.sync();\n  metodoValue.items[0].insertText(\"Provas e Trabalhos\", Word.InsertLocation.replace);\n  await context.sync();\n\n  // Bibliografia paragraph -> gets the old Docente text.\n  paragraphs.items[15].insertText(\n    \"5840560 - Marco Antonio Carvalho Pereira\",\n    Word.InsertLocation.replace\n  );\n\n  await context.sync();\n};\n", "ps1": "# The document content got reshuffled: each \"slot\" paragraph below now\n# displays the text that used to belong to a different slot. We rewrite\n# each paragraph's text in place (line breaks encoded as Chr(11) so they\n# serialize back out as <w:br/> between runs, exactly like before).\n\n$d = $word.ActiveDocument\n$BR = [char]11\n\n# Objetivos paragraph -> gets the (short) \"Programa resumido\" list.\n$d.Paragraphs(6).Range.Text = (\n    \"1. Fundamentos da Gest\u00e3o de Produ\u00e7\u00e3o\" + $BR +\n    \"2. Vis\u00e3o estrat\u00e9gica da Produ\u00e7\u00e3o.\" + $BR +\n    \"3. Projeto em Gest\u00e3o da Produ\u00e7\u00e3o.\" + $BR +\n    \"4. Planejamento e Controle da Produ\u00e7\u00e3o\"\n)\n\n# Docente(s) Respons\u00e1vel(eis) paragraph -> gets the old Objetivos text.\n$d.Paragraphs(8).Range.Text = \"Apresentar conceitos fundamentais de Administra\u00e7\u00e3o da Produ\u00e7\u00e3o.\"\n\n# Programa resumido paragraph -> gets the (long) \"Programa\" list.\n$d.Paragraphs(10).Range.Text = (\n    \"1 - Fundamentos da gest\u00e3o de produ\u00e7\u00e3o: modelo de transforma\u00e7\u00e3o: inputs, processo de transforma\u00e7\u00e3o e outputs. Tipos de Processo de Produ\u00e7\u00e3o\" + $BR +\n    \"2 - Vis\u00e3o estrat\u00e9gica de produ\u00e7\u00e3o: Papel da fun\u00e7\u00e3o produ\u00e7\u00e3o. Objetivos de Desempenho. Estrat\u00e9gias de Produ\u00e7\u00e3o. Ciclo de Vida Produto/Servi\u00e7o.\" + $BR +\n    \"3 \u2013 Projeto em Gest\u00e3o da Produ\u00e7\u00e3o: Tipos de Processos. Projeto de Produtos e Servi\u00e7os. Projeto de Rede de Opera\u00e7\u00f5es Produtivas. Arranjo F\u00edsico.\" + $BR +\n    \"4 - Planejamento e Controle da Produ\u00e7\u00e3o: Material Requirement Planning (MRP), Manufacturing Resources Planning (MPRII), Enterprise Planning (ERP). Produ\u00e7\u00e3o Enxuta. Kanban. Just in Time.\"\n)\n\n# Programa paragraph -> gets the old \"M\u00e9todo\" text.\n$d.Paragraphs(12).Range.Text = \"O desenvolvimento da disciplina ser\u00e1 baseado em leituras, aula expositiva, discuss\u00e3o e resolu\u00e7\u00e3o de estudos de caso.\"\n\n# Avalia\u00e7\u00e3o paragraph keeps its bold \"M\u00e9todo:\"/\"Crit\u00e9rio:\"/\"Norma de\n# recupera\u00e7\u00e3o:\" labels, but the values after each label shift forward.\n# Replace back-to-front so a freshly-inserted value can never collide\n# with (and get matched instead of) a not-yet-processed original value.\n$avaliacao = $d.Paragraphs(14).Range\n\n$rNorma = $avaliacao.Duplicate\n$rNorma.Find.Execute(\n    \"Prova \u00fanica com nota maior ou igual a 5,0 (cinco).\",\n    $false, $true, $false, $false, $false, $true, 1, $false,\n    (\n        \"SLACK, N. et al. Administra\u00e7\u00e3o da Produ\u00e7\u00e3o. 3 ed. S\u00e3o Paulo: Atlas, 2009.\" + $BR +\n        \"CHASE, R. B. E JACOBS, F.R. Administra\u00e7\u00e3o da Produ\u00e7\u00e3o e de Opera\u00e7\u00f5es. 1 ed. Porto Alegre. Bookman. 2009.\" + $BR +\n        \"CORREA, H.L.; CORREA, C.A. Administra\u00e7\u00e3o da Produ\u00e7\u00e3o e Opera\u00e7\u00f5es. 2 ed. S\u00e3o Paulo. Atlas. 2006\"\n    ),\n    2\n)\n\n$rCriterio = $avaliacao.Duplicate\n$rCriterio.Find.Execute(\n    \"Provas e Trabalhos\",\n    $false, $true, $false, $false, $false, $true, 1, $false,\n    \"Prova \u00fanica com nota maior ou igual a 5,0 (cinco).\",\n    2\n)\n\n$rMetodo = $avaliacao.Duplicate\n$rMetodo.Find.Execute(\n    \"O desenvolvimento da disciplina ser\u00e1 baseado em leituras, aula expositiva, discuss\u00e3o e resolu\u00e7\u00e3o de estudos de caso.\",\n    $false, $true, $false, $false, $false, $true, 1, $false,\n    \"Provas e Trabalhos\",\n    2\n)\n\n# Bibliografia paragraph -> gets the old Docente text.\n$d.Paragraphs(16).Range.Text = \"5840560 - Marco Antonio Carvalho Pereira\"\n"}
